$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New testscript row for WAT42 (appended right after the existing WAT41 row, row 28).
$ws.Range("A29").Value = "WAT42"
$ws.Range("B29").Value = "WAT-194"
$ws.Range("C29").Value = 'Verify that "Select an organization where this author has published." is mentioned on top of org list dropdown'
$ws.Range("D29").Value = "Y"

# Match the thin-border look used by every other data row (A:E) in the table.
$ws.Range("A29:E29").Borders.Color = 0
$ws.Range("A29:E29").Borders.LineStyle = 1

# The long description column wraps, like the rest of column C.
$ws.Range("C29").WrapText = $true

# Leave the cursor on the newly-added description cell, as in the authored edit.
[void]$ws.Range("C29").Select()
